# Insert a new blank column before column D, shifting the existing
# "Folder" / "Expected Result..." columns (D:F) one place to the right
# (to E:G), mirroring a user selecting the whole column D and using
# Insert to push the data over.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Range("D1")
$col.EntireColumn.Insert()

# Reflect the resulting selection (whole column D, now empty) as shown
# in the saved file.
$ws.Range("D1:D1048576").Select()
